# Apply the "Added RomaniaFC,SlovakiaFC test data" commit:
#  - Add two new worksheets (Romania, Slovakia) after Spain, cloned from the
#    Belgium sheet (it already has the FCZ.../Dect-Fault row ordering that the
#    new country sheets use) and re-point their header cells.
#  - Tidy up the selection left behind on the Belgium sheet and move the
#    "active" sheet flag off Italy and onto the new Slovakia sheet.

$wb = $excel.ActiveWorkbook

# --- Belgium: clear the old A8:A18 selection to a "select all" state -------
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Activate()
$belgium.Cells.Select()

# --- Romania: clone Belgium (same row layout), rename, retarget header -----
$belgium.Copy($null, $wb.Worksheets.Item("Spain"))
$romania = $wb.Worksheets.Item($wb.Worksheets.Item("Spain").Index + 1)
$romania.Name = "Romania"
$romania.Range("B2").Value = "Romania Market"
$romania.Range("B4").Value = "NGC-4307/T3533/T3549"
$romania.Activate()
$romania.Range("A13").Select()

# --- Slovakia: clone Belgium again, rename, retarget header ----------------
$belgium.Copy($null, $romania)
$slovakia = $wb.Worksheets.Item($romania.Index + 1)
$slovakia.Name = "Slovakia"
$slovakia.Range("B4").Value = "NGC-4306/T3556/T3566"
$slovakia.Range("B2").Value = "Slovakia Market"

# --- Italy loses tabSelected, select-all like Belgium; Slovakia becomes the
#     newly active / last-selected sheet ------------------------------------
$italy = $wb.Worksheets.Item("Italy")
$italy.Activate()
$italy.Cells.Select()

$slovakia.Activate()
$slovakia.Range("B2:B4").Select()
